# Update "paises.xlsx" (Pais sheet) with refreshed COVID-19 country stats
# and fix the three country rows whose rank (and therefore shared string
# slot) moved because the sheet is kept sorted by "Casos totales" desc:
#   Bahamas overtakes Andorra/Chad (rows 155-157)
#   Reunion overtakes San Marino (rows 163-164)
#   Montserrat overtakes Islas Malvinas (rows 213-214)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 01:47"

# Refreshed case counts (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes)
$ws.Range("B4").Value = 5303672
$ws.Range("C4").Value = 52234
$ws.Range("D4").Value = 2749691
$ws.Range("E4").Value = 2386432
$ws.Range("G4").Value = 1357
$ws.Range("H4").Value = 167549
$ws.Range("B5").Value = 3112393
$ws.Range("C5").Value = 54923
$ws.Range("D5").Value = 2243124
$ws.Range("E5").Value = 766170
$ws.Range("G5").Value = 1242
$ws.Range("H5").Value = 103099
$ws.Range("B27").Value = 120421
$ws.Range("C27").Value = 289
$ws.Range("D27").Value = 106746
$ws.Range("E27").Value = 4684
$ws.Range("B46").Value = 57966
$ws.Range("C46").Value = 979
$ws.Range("D46").Value = 46442
$ws.Range("E46").Value = 9291
$ws.Range("G46").Value = 11
$ws.Range("H46").Value = 2233
$ws.Range("B50").Value = 48928
$ws.Range("C50").Value = 938
$ws.Range("D50").Value = 33975
$ws.Range("E50").Value = 13901
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 1052
$ws.Range("B52").Value = 47290
$ws.Range("C52").Value = 423
$ws.Range("D52").Value = 33609
$ws.Range("E52").Value = 12725
$ws.Range("G52").Value = 6
$ws.Range("H52").Value = 956
$ws.Range("B74").Value = 18783
$ws.Range("C74").Value = 289
$ws.Range("D74").Value = 13222
$ws.Range("E74").Value = 5170
$ws.Range("B81").Value = 13722
$ws.Range("C81").Value = 210
$ws.Range("D81").Value = 8154
$ws.Range("E81").Value = 5097
$ws.Range("G81").Value = 12
$ws.Range("H81").Value = 471
$ws.Range("B89").Value = 8360
$ws.Range("C89").Value = 36
$ws.Range("D89").Value = 7632
$ws.Range("E89").Value = 679
$ws.Range("B105").Value = 5223
$ws.Range("C105").Value = 66
$ws.Range("D105").Value = 2849
$ws.Range("E105").Value = 2354
$ws.Range("B112").Value = 4115
$ws.Range("C112").Value = 213
$ws.Range("E112").Value = 1074
$ws.Range("G112").Value = 5
$ws.Range("H112").Value = 128
$ws.Range("B113").Value = 3748
$ws.Range("C113").Value = 52
$ws.Range("D113").Value = 2558
$ws.Range("E113").Value = 1119
$ws.Range("G113").Value = 3
$ws.Range("H113").Value = 71
$ws.Range("B114").Value = 3745
$ws.Range("C114").Value = 81
$ws.Range("D114").Value = 1625
$ws.Range("E114").Value = 2060
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 60
$ws.Range("B125").Value = 2559
$ws.Range("C125").Value = 70
$ws.Range("D125").Value = 1712
$ws.Range("E125").Value = 808
$ws.Range("G125").Value = 9
$ws.Range("H125").Value = 39
$ws.Range("B140").Value = 1385
$ws.Range("C140").Value = 21
$ws.Range("D140").Value = 1157
$ws.Range("E140").Value = 191
$ws.Range("D150").Value = 1065
$ws.Range("E150").Value = 24
$ws.Range("B165").Value = 602
$ws.Range("C165").Value = 34
$ws.Range("E165").Value = 391
$ws.Range("D169").Value = 379
$ws.Range("E169").Value = 13

# Bahamas moves ahead of Andorra and Chad (rows keep their rank order,
# the country names + stats rotate down one row)
$ws.Range("A155").Value = "Bahamas"
$ws.Range("B155").Value = 989
$ws.Range("C155").Value = 44
$ws.Range("D155").Value = 116
$ws.Range("E155").Value = 858
$ws.Range("H155").Value = 15
$ws.Range("A156").Value = "Principado de Andorra"
$ws.Range("B156").Value = 963
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 839
$ws.Range("E156").Value = 72
$ws.Range("H156").Value = 52
$ws.Range("A157").Value = "Republica del Chad"
$ws.Range("B157").Value = 946
$ws.Range("C157").Value = 1
$ws.Range("D157").Value = 859
$ws.Range("E157").Value = 11
$ws.Range("H157").Value = 76

# Reunion overtakes San Marino
$ws.Range("A163").Value = "Reunion"
$ws.Range("B163").Value = 702
$ws.Range("C163").Value = 12
$ws.Range("D163").Value = 631
$ws.Range("E163").Value = 66
$ws.Range("H163").Value = 5
$ws.Range("A164").Value = "San Marino"
$ws.Range("B164").Value = 699
$ws.Range("D164").Value = 657
$ws.Range("E164").Value = 0
$ws.Range("H164").Value = 42

# Montserrat overtakes Islas Malvinas
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
